# Apply edit: rewrite body content to match target revision.
# The document is restructured into two numbered "Scenario" examples of
# IF/FOR template markup (reproducing issue #154 from the upstream repo).
$d = $word.ActiveDocument

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Scenario 1</w:t></w:r></w:p><w:p><w:r><w:t>+++</w:t></w:r><w:r><w:t xml:space="preserve">IF </w:t></w:r><w:r><w:t>list</w:t></w:r><w:r><w:t>+++</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>+++</w:t></w:r><w:r><w:t xml:space="preserve">FOR data IN </w:t></w:r><w:r><w:t>list</w:t></w:r><w:r><w:t>+++ +++</w:t></w:r><w:r><w:t xml:space="preserve">INS </w:t></w:r><w:r><w:t>$</w:t></w:r><w:r><w:t>data</w:t></w:r><w:r><w:t>+++</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>+++</w:t></w:r><w:r><w:t>END-FOR</w:t></w:r><w:r><w:t xml:space="preserve"> data</w:t></w:r><w:r><w:t>+++ +++</w:t></w:r><w:r><w:t>END-IF</w:t></w:r><w:r><w:t>+++</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Scenario 2</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">+++IF list+++ +++IF </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>list[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>0]+++ +++END-IF+++</w:t></w:r></w:p><w:p><w:r><w:t>+++END-IF+++</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
